$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Mirror the Control-table formatting (A2:E9) onto a second block (H2:L9)
#    so the new "HumMod" table has identical borders/fonts/alignment to the
#    existing "QCP" table.
# ---------------------------------------------------------------------------
$ws.Range("A2:E9").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Row labels / headers that reuse strings already in the shared-string
#    table (order does not matter for these - they are not new strings).
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Asphyxia Test"

$ws.Range("H2").Value = "Time"
$ws.Range("I2").Value = "Control"
$ws.Range("J2").Value = "30 Sec"
$ws.Range("K2").Value = "1 Min"
$ws.Range("L2").Value = "5 Min"

$ws.Range("H3").Value = "Blood Pressure(mmHg)"
$ws.Range("H4").Value = "Cardiac Output(mL/min)"
$ws.Range("H5").Value = "Heart Rate(bpm)"
$ws.Range("H6").Value = "Stroke Volume(mL)"
$ws.Range("H7").Value = "Arterial pO2(mmHg)"
$ws.Range("H8").Value = "Blood pH(unitless)"
$ws.Range("I8").Value = "7.37/7.37"
$ws.Range("H9").Value = "Arterial/Venous"

# ---------------------------------------------------------------------------
# 3. New strings - must be entered in this exact order so the shared-string
#    table grows the same way it did in the authored edit (HumMod, QCP,
#    120/79, then the footnote).
# ---------------------------------------------------------------------------
$ws.Range("J1").Value = "HumMod"
$ws.Range("C1").Value = "QCP"
$ws.Range("I3").Value = "120/79"

# ---------------------------------------------------------------------------
# 4. Numeric HumMod results (no shared strings involved).
# ---------------------------------------------------------------------------
$ws.Range("I4").Value = 5468
$ws.Range("I5").Value = 72
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 93

# ---------------------------------------------------------------------------
# 5. Merge the cells in the new table that mirror the merges in the
#    original table (B8:B9, C8:C9, D8:D9, E8:E9 -> I8:I9, J8:J9, K8:K9, L8:L9)
# ---------------------------------------------------------------------------
$ws.Range("I8:I9").Merge()
$ws.Range("J8:J9").Merge()
$ws.Range("K8:K9").Merge()
$ws.Range("L8:L9").Merge()

# ---------------------------------------------------------------------------
# 6. Footnote row (row 11), highlighted in yellow. Build the two new cell
#    styles on scratch cells far away from the used range, copy those
#    formats onto the real targets, and then clear the scratch cells so
#    they leave no trace in the saved sheet.
# ---------------------------------------------------------------------------
$scratchNote = $ws.Range("Z1")
$scratchNote.Font.Name = "Arial"
$scratchNote.Font.Size = 12
$scratchNote.Interior.Color = 65535
$scratchNote.HorizontalAlignment = -4108
$scratchNote.VerticalAlignment = -4160
$scratchNote.WrapText = $true

$scratchFill = $ws.Range("Z2")
$scratchFill.Interior.Color = 65535
$scratchFill.HorizontalAlignment = -4108

$scratchNote.Copy()
$ws.Range("A11").PasteSpecial(-4122)

$scratchFill.Copy()
$ws.Range("B11:K11").PasteSpecial(-4122)

$scratchNote.Clear()
$scratchFill.Clear()

$ws.Range("A11").Value = "The workaround available that allowed this lab to be completed in QCP is not present in HumMod."
$ws.Range("A11:K11").Merge()

# ---------------------------------------------------------------------------
# 7. Column A width + selection, matching the authored workbook.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 9.140625
$ws.Range("F13").Select()

Write-Host "done"
